$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10").Value = 44056
$ws.Range("A10").NumberFormat = "yyyy-mm-dd"
$ws.Range("B10").Value = "NIFTY"
$ws.Range("C10").Value = "SELL"
$ws.Range("D10").Value = 11350.55
$ws.Range("E10").Value = 11339.9
$ws.Range("F10").Value = 11361.2
